$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "SimWell (Laval, QC.)" -- text unchanged, kept as-is (Word's internal
#    proofing split of this run into two runs around a w:proofErr pair is a
#    cosmetic spell-check artifact, not a content edit).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 2) authentication bullet -- text unchanged (same proofing-only artifact).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 3) FGO Gacha Simulator bullet -- text unchanged (same proofing-only
#    artifact) except the tech-stack text further down in the same
#    paragraph also stays the same text.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 4) Script Convenience Store hyperlink display text + URL slug change.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("https://scripts.aaanh.app", $true, $false, $false, $false, $false, $true, 1, $false, "https://script.aaanh.app/", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) GCES Concordia bullet -- text unchanged.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 6) Languages line: "Python, Go" -> "Golang, Python" (reordered + renamed).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Languages: TypeScript, Java, C, C++, Python, Go, Rust, HTML, CSS, PowerShell, bash.", $true, $false, $false, $false, $false, $true, 1, $false, "Languages: TypeScript, Java, C, C++, Golang, Python, Rust, HTML, CSS, PowerShell, bash.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 7) Platforms -> Platform (singular).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Platforms: Linux, Windows, macOS, qemu, Hyper-V, Docker, cri-o, kubernetes, helm, Ansible.", $true, $false, $false, $false, $false, $true, 1, $false, "Platform: Linux, Windows, macOS, qemu, Hyper-V, Docker, cri-o, kubernetes, helm, Ansible.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 8) Networking line -- text unchanged.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 9) Security line -- text unchanged.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 10) Database -> Data, fix CodmosDB -> CosmosDB, append more datastores.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Database: postgres, mysql, CodmosDB, mariadb, redis, memcache.", $true, $false, $false, $false, $false, $true, 1, $false, "Data: postgres, mysql, CosmosDB, mariadb, redis, memcache, Cassandra, Kafka, RabbitMQ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 11) New "Linguistics" bullet added right after the Data bullet, same list
#     (numId 10) and run formatting as its neighbours.
# ---------------------------------------------------------------------------
$dataPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "Data: postgres*") {
        $dataPara = $cand
        break
    }
}
$dataPara.Range.InsertParagraphAfter()
$newIndex = $dataPara.Index + 1
$linguisticsPara = $d.Paragraphs($newIndex)
$linguisticsPara.Range.Text = "Linguistics: Vietnamese (native), English (native), French (B1), Japanese (JLPT N5), German (A1)"

# ---------------------------------------------------------------------------
# 12) HackConcordia hyperlink -- text unchanged.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 13) Remove the "DEC in Applied Sciences ..." education bullet entirely.
# ---------------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "DEC in Applied Sciences*") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 14) Add the "FollowedHyperlink" character style (mirrors the built-in
#     Hyperlink style already present in the template).
# ---------------------------------------------------------------------------
$followed = $d.Styles.Add("FollowedHyperlink", 2)
$followed.BaseStyle = $d.Styles("DefaultParagraphFont")
$followed.Priority = 99
$followed.UnhideWhenUsed = $true
$followed.Font.Underline = 1
$followed.Font.Color = 7491477
Write-Output "Edit script completed."
